$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force target cells to remain text (not auto-converted to numbers/dates)
# by setting the number format to Text before assigning values.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.003.02"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -1.14%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.820.86"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -0.34%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.002"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.31%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "310.83"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -1.16%  "
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -0.23%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4491"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +5.01%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3691"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -0.12%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07305"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +0.63%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.8559"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -1.07%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "20.74"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -1.67%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.818.62"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -0.28%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "6.633"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -1.21%  "
$ws.Range("B14").NumberFormat = "@"
$ws.Range("B14").Value = "TRON"
$ws.Range("C14").NumberFormat = "@"
$ws.Range("C14").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.07117"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +0.16%  "
$ws.Range("B15").NumberFormat = "@"
$ws.Range("B15").Value = "Polkadot"
$ws.Range("C15").NumberFormat = "@"
$ws.Range("C15").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "5.325"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +0.16%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "92.12"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +3.56%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.003"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -0.23%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008771"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -1.15%  "
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -0.35%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "14.95"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -1.03%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "27.037.01"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -1.07%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.165"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +0.27%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "10.91"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +0.51%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.989"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -0.90%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "151.93"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -0.67%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.223"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +3.27%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.43"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -0.14%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "5.239"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -0.31%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "116.38"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -0.31%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.08861"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -0.53%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.7533"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -0.87%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.180"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -1.89%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.957"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +4.09%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.449"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -0.26%  "
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -0.33%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.093"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -1.87%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.01965"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -0.93%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.05229"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -1.17%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.5311"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +4.97%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.891"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +0.56%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "7.120"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -1.20%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1703"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +0.22%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.5227"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +9.69%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "8.492"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -2.40%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "10.65"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +0.07%  "
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +7.21%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "105.44"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -2.17%  "
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -0.33%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.666"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -0.25%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.06379"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +0.08%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.9195"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -0.06%  "
